$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "G1"
$ws.Range("B66").Value = "Test1"
$ws.Range("C66").Value = 45893
$ws.Range("C66").NumberFormat = $ws.Range("C65").NumberFormat
$ws.Range("D66").Value = 0.7345771463238852
$ws.Range("E66").Value = 0
$ws.Range("F66").Value = -0.01

$ws.Range("A67").Value = "G2"
$ws.Range("B67").Value = "sedrftgyhuioygtfrd"
$ws.Range("C67").Value = 45893
$ws.Range("C67").NumberFormat = $ws.Range("C65").NumberFormat
$ws.Range("D67").Value = 0.7345771463238852
$ws.Range("E67").Value = 0
$ws.Range("F67").Value = -0.01
